# Switch license from BY-NC to BY-SA
# Also tidies a couple of split runs on the title slide and nudges the
# license textbox's horizontal offset, mirroring the authored commit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 (title slide): merge a few runs that had been split in two.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

# "Blue Waters Petascale" + " Semester Curriculum v1.0"
#   -> "Blue Waters Petascale Semester Curriculum v1.0"
$full1 = $titleRange.Text
$target1 = "Blue Waters Petascale Semester Curriculum v1.0"
$start1 = $full1.IndexOf("Blue Waters Petascale") + 1
$titleRange.Characters($start1, $target1.Length).Text = $target1

# "Unit " + "10: Productivity and " -> "Unit 10: Productivity and "
$full1 = $titleRange.Text
$target2 = "Unit 10: Productivity and "
$start2 = $full1.IndexOf("Unit ") + 1
$titleRange.Characters($start2, $target2.Length).Text = $target2

# "by " + "Michael N. " -> "by Michael N. "
$full1 = $titleRange.Text
$target3 = "by Michael N. "
$start3 = $full1.IndexOf("by Michael N. ") + 1
$titleRange.Characters($start3, $target3.Length).Text = $target3

# ---------------------------------------------------------------------
# Slide 2 (license slide): BY-NC -> BY-SA, in both the visible license
# name and the URL that follows it. Nudge the textbox's x-offset.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$licenseShape = $slide2.Shapes.Item(1)
$licenseRange = $licenseShape.TextFrame.TextRange

# "CC BY-NC 4.0. ..." -> "CC BY-SA 4.0. ..."
$full2 = $licenseRange.Text
$ncStart = $full2.IndexOf("BY-NC ") + 1
$licenseRange.Characters($ncStart, "BY-NC ".Length).Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0"
#   -> "https://creativecommons.org/licenses/by-sa/4.0"
$full2 = $licenseRange.Text
$urlSuffixStart = $full2.IndexOf("creativecommons.org/licenses/by-nc/4.0") + 1
$newSuffix = "creativecommons.org/licenses/by-sa/4.0"
$licenseRange.Characters($urlSuffixStart, "creativecommons.org/licenses/by-nc/4.0".Length).Text = $newSuffix

# Shape x-offset: 566059 EMU -> 566057 EMU (914400 EMU per inch, 12700 EMU per point)
$licenseShape.Left = 566057 / 12700
